# The post that used to live in row 756 ("「象の群れ。ドローンの音から逃げている」 ...")
# was removed from the spreadsheet. Delete that entire row; Excel will
# automatically shift every following row (757-828) up by one, which is
# exactly the renumbering shown in the diff (dimension becomes A1:C827).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(756).Delete()
